{"js": "const replacements = [\n  [\"392\u00f73=130, 2\", \"281\u00f72=140, 1\"],\n  [\"855\u00f77=122, 1\", \"583\u00f72=291, 1\"],\n  [\"974\u00f75=194, 4\", \"236\u00f72=118, 0\"],\n  [\"512\u00f78=64, 0\", \"312\u00f74=78, 0\"],\n  [\"166\u00f73=55, 1\", \"368\u00f75=73, 3\"],\n  [\"122\u00f74=30, 2\", \"273\u00f73=91, 0\"],\n  [\"672\u00f74=168, 0\", \"966\u00f79=107, 3\"],\n  [\"975\u00f77=139, 2\", \"742\u00f79=82, 4\"],\n  [\"694\u00f73=231, 1\", \"173\u00f73=57, 2\"],\n  [\"701\u00f77=100, 1\", \"256\u00f75=51, 1\"],\n  [\"372\u00f78=46, 4\", \"534\u00f74=133, 2\"],\n  [\"903\u00f77=129, 0\", \"227\u00f78=28, 3\"],\n  [\"564\u00f75=112, 4\", \"342\u00f78=42, 6\"],\n  [\"988\u00f78=123, 4\", \"412\u00f78=51, 4\"],\n  [\"250\u00f75=50, 0\", \"540\u00f79=60, 0\"],\n  [\"524\u00f73=174, 2\", \"976\u00f74=244, 0\"],\n  [\"362\u00f78=45, 2\", \"646\u00f72=323, 0\"],\n  [\"309\u00f72=154, 1\", \"960\u00f77=137, 1\"],\n  [\"723\u00f76=120, 3\", \"312\u00f74=78, 0\"],\n  [\"514\u00f78=64, 2\", \"915\u00f78=114, 3\"],\n  [\"163\u00f73=54, 1\", \"460\u00f76=76, 4\"],\n  [\"905\u00f76=150, 5\", \"802\u00f77=114, 4\"],\n  [\"696\u00f75=139, 1\", \"105\u00f79=11, 6\"],\n  [\"918\u00f78=114, 6\", \"363\u00f78=45, 3\"],\n  [\"447\u00f78=55, 7\", \"744\u00f75=148, 4\"]\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"392\u00f73=130, 2\", \"281\u00f72=140, 1\"),\n    @(\"855\u00f77=122, 1\", \"583\u00f72=291, 1\"),\n    @(\"974\u00f75=194, 4\", \"236\u00f72=118, 0\"),\n    @(\"512\u00f78=64, 0\", \"312\u00f74=78, 0\"),\n    @(\"166\u00f73=55, 1\", \"368\u00f75=73, 3\"),\n    @(\"122\u00f74=30, 2\", \"273\u00f73=91, 0\"),\n    @(\"672\u00f74=168, 0\", \"966\u00f79=107, 3\"),\n    @(\"975\u00f77=139, 2\", \"742\u00f79=82, 4\"),\n    @(\"694\u00f73=231, 1\", \"173\u00f73=57, 2\"),\n    @(\"701\u00f77=100, 1\", \"256\u00f75=51, 1\"),\n    @(\"372\u00f78=46, 4\", \"534\u00f74=133, 2\"),\n    @(\"903\u00f77=129, 0\", \"227\u00f78=28, 3\"),\n    @(\"564\u00f75=112, 4\", \"342\u00f78=42, 6\"),\n    @(\"988\u00f78=123, 4\", \"412\u00f78=51, 4\"),\n    @(\"250\u00f75=50, 0\", \"540\u00f79=60, 0\"),\n    @(\"524\u00f73=174, 2\", \"976\u00f74=244, 0\"),\n    @(\"362\u00f78=45, 2\", \"646\u00f72=323, 0\"),\n    @(\"309\u00f72=154, 1\", \"960\u00f77=137, 1\"),\n    @(\"723\u00f76=120, 3\", \"312\u00f74=78, 0\"),\n    @(\"514\u00f78=64, 2\", \"915\u00f78=114, 3\"),\n    @(\"163\u00f73=54, 1\", \"460\u00f76=76, 4\"),\n    @(\"905\u00f76=150, 5\", \"802\u00f77=114, 4\"),\n    @(\"696\u00f75=139, 1\", \"105\u00f79=11, 6\"),\n    @(\"918\u00f78=114, 6\", \"363\u00f78=45, 3\"),\n    @(\"447\u00f78=55, 7\", \"744\u00f75=148, 4\"),\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
